# Update cryptocurrency price/volume snapshot (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.871.36'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.812.98'
$ws.Range('E3').Value = '  +0.91%  '

$ws.Range('E4').Value = '  +0.32%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.86'
$ws.Range('E5').Value = '  +0.19%  '

$ws.Range('E6').Value = '  +0.25%  '

$ws.Range('E7').Value = '  +1.84%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3686'
$ws.Range('E8').Value = '  -0.72%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07363'
$ws.Range('E9').Value = '  +1.69%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8689'
$ws.Range('E10').Value = '  +1.52%  '

$ws.Range('E11').Value = '  -0.14%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.781.44'
$ws.Range('E12').Value = '  -0.90%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.344'
$ws.Range('E13').Value = '  +0.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07073'
$ws.Range('E14').Value = '  +0.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.80'
$ws.Range('E15').Value = '  +1.71%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.488'
$ws.Range('E16').Value = '  -0.27%  '

$ws.Range('E17').Value = '  +0.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008690'
$ws.Range('E18').Value = '  +0.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.74'
$ws.Range('E20').Value = '  +0.77%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.932.37'
$ws.Range('E21').Value = '  +0.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.338'
$ws.Range('E22').Value = '  +0.91%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.54'
$ws.Range('E23').Value = '  -0.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.053.19'
$ws.Range('E24').Value = '  +1.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.900'
$ws.Range('E25').Value = '  -0.42%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.11'
$ws.Range('E26').Value = '  +1.07%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.168'
$ws.Range('E27').Value = '  +0.35%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.33'
$ws.Range('E28').Value = '  +0.83%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.312'
$ws.Range('E29').Value = '  +2.07%  '

$ws.Range('E30').Value = '  +1.29%  '

$ws.Range('E31').Value = '  +0.98%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7655'
$ws.Range('E32').Value = '  +0.90%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.158'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.497'
$ws.Range('E34').Value = '  +1.15%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.905'
$ws.Range('E35').Value = '  +0.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.002'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.087'
$ws.Range('E37').Value = '  -2.34%  '

$ws.Range('E38').Value = '  +0.89%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05284'
$ws.Range('E39').Value = '  +1.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.940'
$ws.Range('E40').Value = '  +1.60%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.269'
$ws.Range('E41').Value = '  +1.91%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5305'
$ws.Range('E42').Value = '  +1.26%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.349'
$ws.Range('E43').Value = '  -1.61%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1662'
$ws.Range('E44').Value = '  +0.92%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.404'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4925'
$ws.Range('E46').Value = '  -2.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.39'
$ws.Range('E47').Value = '  +1.64%  '

$ws.Range('E48').Value = '  +0.27%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.64'
$ws.Range('E49').Value = '  -0.61%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.667'
$ws.Range('E50').Value = '  +1.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06277'
$ws.Range('E51').Value = '  -0.10%  '
